$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Price" column (D) holds strings that must stay literal text
# (preserving exact digits/trailing zeros, e.g. "1.00", "0.998").
# Whenever the new value would otherwise parse as a pure number,
# force the cell to Text format first so Excel keeps it as a string
# instead of silently coercing it into a numeric value.

$ws.Range('D4').NumberFormat = '@'
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D51').NumberFormat = '@'

$ws.Range('D2').Value = '54.008.66'
$ws.Range('E2').Value = '  -9.01%  '

$ws.Range('D3').Value = '2.391.44'
$ws.Range('E3').Value = '  -16.14%  '

$ws.Range('D4').Value = '0.998'
$ws.Range('E4').Value = '  -0.22%  '

$ws.Range('D5').Value = '459.91'
$ws.Range('E5').Value = '  -8.60%  '

$ws.Range('D6').Value = '130.49'
$ws.Range('E6').Value = '  -4.47%  '

$ws.Range('D7').Value = '0.994'
$ws.Range('E7').Value = '  -0.74%  '

$ws.Range('D8').Value = '0.483'
$ws.Range('E8').Value = '  -8.41%  '

$ws.Range('D9').Value = '2.404.40'
$ws.Range('E9').Value = '  -15.72%  '

$ws.Range('D10').Value = '0.0942'
$ws.Range('E10').Value = '  -8.80%  '

$ws.Range('D11').Value = '5.32'
$ws.Range('E11').Value = '  -11.40%  '

$ws.Range('E12').Value = '  -8.66%  '

$ws.Range('E13').Value = '  -4.27%  '

$ws.Range('D14').Value = '2.787.03'
$ws.Range('E14').Value = '  -16.96%  '

$ws.Range('D15').Value = '53.803.23'
$ws.Range('E15').Value = '  -9.78%  '

$ws.Range('D16').Value = '19.52'
$ws.Range('E16').Value = '  -10.14%  '

$ws.Range('D17').Value = '0.0000128'
$ws.Range('E17').Value = '  -5.75%  '

$ws.Range('D18').Value = '2.392.83'
$ws.Range('E18').Value = '  -16.64%  '

$ws.Range('E19').Value = '  -11.68%  '

$ws.Range('D20').Value = '309.49'
$ws.Range('E20').Value = '  -11.62%  '

$ws.Range('D21').Value = '9.28'
$ws.Range('E21').Value = '  -16.63%  '

$ws.Range('D22').Value = '0.999'
$ws.Range('E22').Value = '  +0.20%  '

$ws.Range('E23').Value = '  +0.24%  '

$ws.Range('D24').Value = '5.37'
$ws.Range('E24').Value = '  -14.35%  '

$ws.Range('D25').Value = '56.14'
$ws.Range('E25').Value = '  -11.02%  '

$ws.Range('D26').Value = '1.00'
$ws.Range('E26').Value = '  +0.38%  '

$ws.Range('E27').Value = '  -12.02%  '

$ws.Range('E28').Value = '  -11.82%  '

$ws.Range('D29').Value = '2.474.29'
$ws.Range('E29').Value = '  -17.17%  '

$ws.Range('D30').Value = '7.08'
$ws.Range('E30').Value = '  -5.79%  '

$ws.Range('E31').Value = '  -0.36%  '

$ws.Range('D32').Value = '0.0₃0706'
$ws.Range('E32').Value = '  -14.19%  '

$ws.Range('D33').Value = '149.10'
$ws.Range('E33').Value = '  -1.27%  '

$ws.Range('D34').Value = '17.59'
$ws.Range('E34').Value = '  -7.64%  '

$ws.Range('E35').Value = '  -13.70%  '

$ws.Range('D36').Value = '4.99'
$ws.Range('E36').Value = '  -7.25%  '

$ws.Range('D37').Value = '3.44'
$ws.Range('E37').Value = '  -18.00%  '

$ws.Range('E38').Value = '  -9.63%  '

$ws.Range('E39').Value = '  -16.00%  '

$ws.Range('D40').Value = '33.74'
$ws.Range('E40').Value = '  -7.88%  '

$ws.Range('D41').Value = '0.989'
$ws.Range('E41').Value = '  -1.07%  '

$ws.Range('D42').Value = '0.603'
$ws.Range('E42').Value = '  -4.72%  '

$ws.Range('D43').Value = '3.29'
$ws.Range('E43').Value = '  -6.95%  '

$ws.Range('E44').Value = '  -6.36%  '

$ws.Range('D45').Value = '10.15'
$ws.Range('E45').Value = '  -1.97%  '

$ws.Range('B46').Value = 'Stacks'
$ws.Range('C46').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D46').Value = '1.23'
$ws.Range('E46').Value = '  -11.38%  '

$ws.Range('B47').Value = 'Maker'
$ws.Range('C47').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D47').Value = '1.946.20'
$ws.Range('E47').Value = '  -13.08%  '

$ws.Range('D48').Value = '0.0217'
$ws.Range('E48').Value = '  -4.36%  '

$ws.Range('D49').Value = '0.0864'
$ws.Range('E49').Value = '  -2.88%  '

$ws.Range('D50').Value = '4.26'
$ws.Range('E50').Value = '  -6.91%  '

$ws.Range('D51').Value = '16.40'
$ws.Range('E51').Value = '  -16.58%  '
